$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to text format first so purely numeric-looking
# strings like "571.61" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.694.37"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.606.33"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "571.61"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "142.69"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "2.627.60"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "0.368"
$ws.Range("E13").Value = "  +6.72%  "
$ws.Range("D14").Value = "3.074.51"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "60.721.43"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "23.57"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "2.624.20"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "4.69"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").Value = "  +9.29%  "
$ws.Range("D21").Value = "348.42"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("D22").Value = "7.06"
$ws.Range("E22").Value = "  +13.53%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "0.517"
$ws.Range("E24").Value = "  +13.40%  "
$ws.Range("D25").Value = "63.79"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  +5.86%  "
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  +10.61%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "161.53"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "19.51"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  +4.95%  "
$ws.Range("D36").Value = "0.958"
$ws.Range("E36").Value = "  +8.55%  "
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("E38").Value = "  +6.39%  "
$ws.Range("D39").Value = "37.66"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "0.857"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").Value = "296.89"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "138.69"
$ws.Range("E43").Value = "  +10.29%  "
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "0.0985"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "0.0551"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").Value = "0.605"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").Value = "10.70"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "19.62"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("D51").Value = "4.81"
$ws.Range("E51").Value = "  +6.31%  "

# Restore the original cell style (no explicit style index) now that
# the values are stored as text, matching the workbook's original formatting.
$ws.Range("D2:D51").Style = "Normal"
